$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("2023-12-11 18:16:37", 0.0006000000000000001),
    @("2023-12-11 18:17:14", 0.0018),
    @("2023-12-11 18:17:32", 0.0008),
    @("2023-12-11 18:17:37", 0.0004),
    @("2023-12-11 18:17:54", 0.0008)
)

$startRow = 201
for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $startRow + $i
    $ws.Cells.Item($row, 1).Value = $data[$i][0]
    $ws.Cells.Item($row, 2).Value = $data[$i][1]
}

$wb.Save()
